$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 32118.031
$ws.Cells.Item(17, 10).Value = 32118.031
$ws.Cells.Item(17, 12).Value = 96354.09299999999
$ws.Cells.Item(17, 14).Value = -96690.09299999999
$ws.Cells.Item(88, 8).Value = 7607084
$ws.Cells.Item(88, 9).Value = 1239
$ws.Cells.Item(88, 10).Value = 10738902
$ws.Cells.Item(88, 11).Value = 1239
$ws.Cells.Item(88, 12).Value = 10738902
$ws.Cells.Item(88, 13).Value = -833
$ws.Cells.Item(88, 14).Value = -10739714
$ws.Cells.Item(91, 8).Value = 7607084
$ws.Cells.Item(91, 9).Value = 1239
$ws.Cells.Item(91, 10).Value = 10738902
$ws.Cells.Item(91, 11).Value = 1239
$ws.Cells.Item(91, 12).Value = 10738902
$ws.Cells.Item(91, 13).Value = 165
$ws.Cells.Item(91, 14).Value = -10741710
$ws.Cells.Item(111, 8).Value = 532.75
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 532.75
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 12).Value = 1598.25
$ws.Cells.Item(111, 13).ClearContents()
$ws.Cells.Item(111, 14).Value = -7732.25
$ws.Cells.Item(112, 8).Value = 12135.889
$ws.Cells.Item(112, 10).Value = 12466.914
$ws.Cells.Item(112, 12).Value = 37400.742
$ws.Cells.Item(112, 14).Value = -39616.742
$ws.Cells.Item(127, 8).Value = 1077.5238
$ws.Cells.Item(127, 9).Value = 449.15384
$ws.Cells.Item(127, 10).Value = 2098.625
$ws.Cells.Item(127, 11).Value = 1347.46152
$ws.Cells.Item(127, 12).Value = 6295.875
$ws.Cells.Item(127, 13).Value = 3612.53848
$ws.Cells.Item(127, 14).Value = -16215.875
$ws.Cells.Item(137, 8).Value = 1105.4546
$ws.Cells.Item(137, 9).Value = 912.8909
$ws.Cells.Item(137, 10).Value = 2068.2727
$ws.Cells.Item(137, 11).Value = 2738.6727
$ws.Cells.Item(137, 12).Value = 6204.8181
$ws.Cells.Item(137, 13).Value = -188.6727000000001
$ws.Cells.Item(137, 14).Value = -11304.8181
$ws.Cells.Item(138, 8).Value = 1484.7084
$ws.Cells.Item(138, 9).Value = 864.54
$ws.Cells.Item(138, 10).Value = 2894.182
$ws.Cells.Item(138, 11).Value = 2593.62
$ws.Cells.Item(138, 12).Value = 8682.545999999998
$ws.Cells.Item(138, 13).Value = 2546.38
$ws.Cells.Item(138, 14).Value = -18962.546
$ws.Cells.Item(141, 8).Value = 2845.1094
$ws.Cells.Item(141, 9).Value = 1114.4043
$ws.Cells.Item(141, 10).Value = 7630
$ws.Cells.Item(141, 11).Value = 3343.2129
$ws.Cells.Item(141, 12).Value = 22890
$ws.Cells.Item(141, 13).Value = 1836.7871
$ws.Cells.Item(141, 14).Value = -33250
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 15000
$ws.Cells.Item(6, 10).Value = 15000
$ws.Cells.Item(6, 12).Value = 15000
$ws.Cells.Item(6, 14).Value = -15346
$ws.Cells.Item(61, 8).Value = 911.0441
$ws.Cells.Item(61, 9).Value = 731.3728599999999
$ws.Cells.Item(61, 10).Value = 2088.889
$ws.Cells.Item(61, 11).Value = 731.3728599999999
$ws.Cells.Item(61, 12).Value = 2088.889
$ws.Cells.Item(61, 13).Value = -519.3728599999999
$ws.Cells.Item(61, 14).Value = -2512.889
$ws.Cells.Item(74, 8).Value = 1047.2258
$ws.Cells.Item(74, 9).Value = 1049.9824
$ws.Cells.Item(74, 10).Value = 1015.8
$ws.Cells.Item(74, 11).Value = 1049.9824
$ws.Cells.Item(74, 12).Value = 1015.8
$ws.Cells.Item(74, 13).Value = -175.9824000000001
$ws.Cells.Item(74, 14).Value = -2763.8
$ws.Cells.Item(77, 8).Value = 1047.2258
$ws.Cells.Item(77, 9).Value = 1049.9824
$ws.Cells.Item(77, 10).Value = 1015.8
$ws.Cells.Item(77, 11).Value = 5249.912
$ws.Cells.Item(77, 12).Value = 5079
$ws.Cells.Item(77, 13).Value = -881.9120000000003
$ws.Cells.Item(77, 14).Value = -13815
$ws.Cells.Item(132, 8).Value = 25826428
$ws.Cells.Item(132, 9).Value = 30304104
$ws.Cells.Item(132, 10).Value = 7356016
$ws.Cells.Item(132, 11).Value = 90912312
$ws.Cells.Item(132, 12).Value = 22068048
$ws.Cells.Item(132, 13).Value = -90909782
$ws.Cells.Item(132, 14).Value = -22073108
$ws.Cells.Item(136, 8).Value = 911.0441
$ws.Cells.Item(136, 9).Value = 731.3728599999999
$ws.Cells.Item(136, 10).Value = 2088.889
$ws.Cells.Item(136, 11).Value = 2194.11858
$ws.Cells.Item(136, 12).Value = 6266.667
$ws.Cells.Item(136, 13).Value = 355.8814200000002
$ws.Cells.Item(136, 14).Value = -11366.667
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1952016.1
$ws.Cells.Item(134, 9).Value = 813.30554
$ws.Cells.Item(134, 10).Value = 5296935.5
$ws.Cells.Item(134, 11).Value = 2439.91662
$ws.Cells.Item(134, 12).Value = 15890806.5
$ws.Cells.Item(134, 13).Value = 95.08338000000003
$ws.Cells.Item(134, 14).Value = -15895876.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(11, 8).Value = 7000
$ws.Cells.Item(11, 9).Value = 7000
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 7000
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 14).ClearContents()
$ws.Cells.Item(11, 13).Value = -6860
$ws.Cells.Item(31, 8).Value = 1233.4266
$ws.Cells.Item(31, 9).Value = 910.7818
$ws.Cells.Item(31, 10).Value = 2120.7
$ws.Cells.Item(31, 11).Value = 910.7818
$ws.Cells.Item(31, 12).Value = 2120.7
$ws.Cells.Item(31, 13).Value = -615.7818
$ws.Cells.Item(31, 14).Value = -2710.7
$ws.Cells.Item(34, 8).Value = 1233.4266
$ws.Cells.Item(34, 9).Value = 910.7818
$ws.Cells.Item(34, 10).Value = 2120.7
$ws.Cells.Item(34, 11).Value = 910.7818
$ws.Cells.Item(34, 12).Value = 2120.7
$ws.Cells.Item(34, 13).Value = -708.7818
$ws.Cells.Item(34, 14).Value = -2524.7
$ws.Cells.Item(58, 8).Value = 17242202
$ws.Cells.Item(58, 9).Value = 25000598
$ws.Cells.Item(58, 10).Value = 1324.6111
$ws.Cells.Item(58, 11).Value = 25000598
$ws.Cells.Item(58, 12).Value = 1324.6111
$ws.Cells.Item(58, 13).Value = -25000395
$ws.Cells.Item(58, 14).Value = -1730.6111
$ws.Cells.Item(132, 8).Value = 8548385
$ws.Cells.Item(132, 9).Value = 1129
$ws.Cells.Item(132, 11).Value = 3387
$ws.Cells.Item(132, 13).Value = -857
$ws.Cells.Item(134, 8).Value = 1306.2051
$ws.Cells.Item(134, 9).Value = 1056.2142
$ws.Cells.Item(134, 11).Value = 3168.6426
$ws.Cells.Item(134, 13).Value = -633.6425999999997
$ws.Cells.Item(136, 8).Value = 17242202
$ws.Cells.Item(136, 9).Value = 25000598
$ws.Cells.Item(136, 10).Value = 1324.6111
$ws.Cells.Item(136, 11).Value = 75001794
$ws.Cells.Item(136, 12).Value = 3973.8333
$ws.Cells.Item(136, 13).Value = -74999244
$ws.Cells.Item(136, 14).Value = -9073.8333
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 275050300
$ws.Cells.Item(9, 9).Value = 100000000
$ws.Cells.Item(9, 10).Value = 333400400
$ws.Cells.Item(9, 11).Value = 300000000
$ws.Cells.Item(9, 12).Value = 1000201200
$ws.Cells.Item(9, 13).Value = -299999776
$ws.Cells.Item(9, 14).Value = -1000201648
$ws.Cells.Item(114, 8).Value = 776.61536
$ws.Cells.Item(114, 9).Value = 210.625
$ws.Cells.Item(114, 10).Value = 1682.2
$ws.Cells.Item(114, 11).Value = 631.875
$ws.Cells.Item(114, 12).Value = 5046.6
$ws.Cells.Item(114, 13).Value = 2622.125
$ws.Cells.Item(114, 14).Value = -11554.6
$ws.Cells.Item(117, 8).Value = 489.2143
$ws.Cells.Item(117, 9).Value = 279.85715
$ws.Cells.Item(117, 10).Value = 698.5714
$ws.Cells.Item(117, 11).Value = 839.5714499999999
$ws.Cells.Item(117, 12).Value = 2095.7142
$ws.Cells.Item(117, 13).Value = 2602.42855
$ws.Cells.Item(117, 14).Value = -8979.7142
$ws.Cells.Item(129, 8).Value = 14493804
$ws.Cells.Item(129, 10).Value = 33334654
$ws.Cells.Item(129, 12).Value = 100003962
$ws.Cells.Item(129, 14).Value = -100013962
$ws.Cells.Item(130, 8).Value = 71430480
$ws.Cells.Item(130, 9).Value = 500000000
$ws.Cells.Item(130, 11).Value = 1500000000
$ws.Cells.Item(130, 13).Value = -1499994980
$ws.Cells.Item(131, 8).Value = 16207988
$ws.Cells.Item(131, 9).Value = 55556790
$ws.Cells.Item(131, 10).Value = 8338226.5
$ws.Cells.Item(131, 11).Value = 166670370
$ws.Cells.Item(131, 12).Value = 25014679.5
$ws.Cells.Item(131, 13).Value = -166665330
$ws.Cells.Item(131, 14).Value = -25024759.5
$ws.Cells.Item(133, 8).Value = 55558730
$ws.Cells.Item(134, 8).Value = 38464492
$ws.Cells.Item(134, 10).Value = 4584
$ws.Cells.Item(134, 12).Value = 13752
$ws.Cells.Item(134, 14).Value = -23892
$ws.Cells.Item(136, 8).Value = 39065936
$ws.Cells.Item(136, 9).Value = 78126230
$ws.Cells.Item(136, 10).Value = 5636.75
$ws.Cells.Item(136, 11).Value = 234378690
$ws.Cells.Item(136, 12).Value = 16910.25
$ws.Cells.Item(136, 13).Value = -234373590
$ws.Cells.Item(136, 14).Value = -27110.25
$ws.Cells.Item(137, 8).Value = 49021104
$ws.Cells.Item(137, 10).Value = 66669044
$ws.Cells.Item(137, 12).Value = 200007132
$ws.Cells.Item(137, 14).Value = -200017332
$ws.Cells.Item(138, 8).Value = 21334762
$ws.Cells.Item(138, 9).Value = 28071194
$ws.Cells.Item(138, 11).Value = 84213582
$ws.Cells.Item(138, 13).Value = -84208442
$ws.Cells.Item(139, 8).Value = 17021896
$ws.Cells.Item(139, 9).Value = 31251066
$ws.Cells.Item(139, 10).Value = 759985.7
$ws.Cells.Item(139, 11).Value = 93753198
$ws.Cells.Item(139, 12).Value = 2279957.1
$ws.Cells.Item(139, 13).Value = -93748058
$ws.Cells.Item(139, 14).Value = -2290237.1
$ws.Cells.Item(140, 8).Value = 19567236
$ws.Cells.Item(140, 9).Value = 28126656
$ws.Cells.Item(140, 10).Value = 2847.4285
$ws.Cells.Item(140, 11).Value = 84379968
$ws.Cells.Item(140, 12).Value = 8542.2855
$ws.Cells.Item(140, 13).Value = -84374788
$ws.Cells.Item(140, 14).Value = -18902.2855
$ws.Cells.Item(141, 8).Value = 41669056
$ws.Cells.Item(141, 9).Value = 50001816
$ws.Cells.Item(141, 11).Value = 150005448
$ws.Cells.Item(141, 13).Value = -150000268
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1446.1578
$ws.Cells.Item(7, 9).Value = 1415.3889
$ws.Cells.Item(7, 10).Value = 2000
$ws.Cells.Item(7, 11).Value = 1415.3889
$ws.Cells.Item(7, 12).Value = 2000
$ws.Cells.Item(7, 13).Value = -1303.3889
$ws.Cells.Item(7, 14).Value = -2224
$ws.Cells.Item(69, 8).Value = 37490.25
$ws.Cells.Item(69, 10).Value = 37490.25
$ws.Cells.Item(69, 12).Value = 37490.25
$ws.Cells.Item(69, 14).Value = -39112.25
$ws.Cells.Item(72, 8).Value = 37490.25
$ws.Cells.Item(72, 10).Value = 37490.25
$ws.Cells.Item(72, 12).Value = 112470.75
$ws.Cells.Item(72, 14).Value = -120582.75
$ws.Cells.Item(126, 8).Value = 1446.1578
$ws.Cells.Item(126, 9).Value = 1415.3889
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 4246.1667
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -1776.1667
$ws.Cells.Item(126, 14).Value = -10940
$ws.Cells.Item(132, 8).Value = 13893292
$ws.Cells.Item(132, 9).Value = 20001370
$ws.Cells.Item(132, 11).Value = 60004110
$ws.Cells.Item(132, 13).Value = -60001580
$ws.Cells.Item(136, 8).Value = 21979732
$ws.Cells.Item(136, 9).Value = 3573369.5
$ws.Cells.Item(136, 11).Value = 10720108.5
$ws.Cells.Item(136, 13).Value = -10717558.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5569526.5
$ws.Cells.Item(132, 9).Value = 16120.529
$ws.Cells.Item(132, 10).Value = 22734600
$ws.Cells.Item(132, 11).Value = 48361.587
$ws.Cells.Item(132, 12).Value = 68203800
$ws.Cells.Item(132, 13).Value = -45831.587
$ws.Cells.Item(132, 14).Value = -68208860
$ws.Cells.Item(136, 8).Value = 9095316
$ws.Cells.Item(136, 9).Value = 11909652
$ws.Cells.Item(136, 10).Value = 2846.8462
$ws.Cells.Item(136, 11).Value = 35728956
$ws.Cells.Item(136, 12).Value = 8540.5386
$ws.Cells.Item(136, 13).Value = -35726406
$ws.Cells.Item(136, 14).Value = -13640.5386
